$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 453.91306
$ws.Range("I33").Value = 483.6842
$ws.Range("J33").Value = 312.5
$ws.Range("K33").Value = 483.6842
$ws.Range("L33").Value = 312.5
$ws.Range("M33").Value = -254.6842
$ws.Range("N33").Value = -770.5

$ws.Range("H49").Value = 684
$ws.Range("I49").Value = 460
$ws.Range("J49").Value = 833.3333
$ws.Range("K49").Value = 1380
$ws.Range("L49").Value = 2499.9999
$ws.Range("M49").Value = -1244
$ws.Range("N49").Value = -2771.9999

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H116").Value = 2505.9375
$ws.Range("I116").Value = 1439
$ws.Range("J116").Value = 2990.9092
$ws.Range("K116").Value = 1439
$ws.Range("L116").Value = 2990.9092
$ws.Range("M116").Value = 2003
$ws.Range("N116").Value = -9874.9092

$ws.Range("H137").Value = 1312.4642
$ws.Range("I137").Value = 785.9722
$ws.Range("J137").Value = 2260.15
$ws.Range("K137").Value = 2357.9166
$ws.Range("L137").Value = 6780.450000000001
$ws.Range("M137").Value = 192.0834
$ws.Range("N137").Value = -11880.45

$ws.Range("H138").Value = 1773.1139
$ws.Range("I138").Value = 1244.6538
$ws.Range("J138").Value = 2032.3585
$ws.Range("K138").Value = 3733.9614
$ws.Range("L138").Value = 6097.0755
$ws.Range("M138").Value = 1406.0386
$ws.Range("N138").Value = -16377.0755

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 21739730
$ws.Range("I2").Value = 31250366
$ws.Range("K2").Value = 31250366
$ws.Range("M2").Value = -31250253

$ws.Range("H32").Value = 17054.967
$ws.Range("I32").Value = 18214.828
$ws.Range("J32").Value = 12544.389
$ws.Range("K32").Value = 18214.828
$ws.Range("L32").Value = 12544.389
$ws.Range("M32").Value = -17927.828
$ws.Range("N32").Value = -13118.389

$ws.Range("H35").Value = 933.3333
$ws.Range("I35").Value = 933.3333
$ws.Range("K35").Value = 933.3333
$ws.Range("M35").Value = -527.3333

$ws.Range("H116").Value = 21739730
$ws.Range("I116").Value = 31250366
$ws.Range("K116").Value = 31250366
$ws.Range("M116").Value = -31248072

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 21739730
$ws.Range("I3").Value = 31250366
$ws.Range("K3").Value = 31250366
$ws.Range("M3").Value = -31250252

$ws.Range("H86").Value = 2900
$ws.Range("I86").Value = 2825
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 2825
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1702
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 2900
$ws.Range("I89").Value = 2825
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 14125
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -8509
$ws.Range("N89").Value = -26232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4068021.5
$ws.Range("I31").Value = 2363.4062
$ws.Range("J31").Value = 18523694
$ws.Range("K31").Value = 2363.4062
$ws.Range("L31").Value = 18523694
$ws.Range("M31").Value = -2068.4062
$ws.Range("N31").Value = -18524284

$ws.Range("H34").Value = 4068021.5
$ws.Range("I34").Value = 2363.4062
$ws.Range("J34").Value = 18523694
$ws.Range("K34").Value = 2363.4062
$ws.Range("L34").Value = 18523694
$ws.Range("M34").Value = -2161.4062
$ws.Range("N34").Value = -18524098

$ws.Range("H86").Value = 333335840
$ws.Range("I86").Value = 500001500
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 500001500
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -500000377
$ws.Range("N86").Value = -6746

$ws.Range("H89").Value = 333335840
$ws.Range("I89").Value = 500001500
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 2500007500
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -2500001884
$ws.Range("N89").Value = -33732

$ws.Range("H107").Value = 584.64813
$ws.Range("I107").Value = 517.1667
$ws.Range("J107").Value = 719.6111
$ws.Range("K107").Value = 517.1667
$ws.Range("L107").Value = 719.6111
$ws.Range("M107").Value = 1402.8333
$ws.Range("N107").Value = -4559.6111

$ws.Range("H132").Value = 3127572.5
$ws.Range("I132").Value = 2112.5454
$ws.Range("J132").Value = 6947579
$ws.Range("K132").Value = 6337.6362
$ws.Range("L132").Value = 20842737
$ws.Range("M132").Value = -3807.6362
$ws.Range("N132").Value = -20847797

$ws.Range("H134").Value = 1237.6
$ws.Range("I134").Value = 1199.75
$ws.Range("K134").Value = 3599.25
$ws.Range("M134").Value = -1064.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1236.0513
$ws.Range("I5").Value = 222.64285
$ws.Range("J5").Value = 1803.56
$ws.Range("K5").Value = 667.9285500000001
$ws.Range("L5").Value = 5410.68
$ws.Range("M5").Value = -555.9285500000001
$ws.Range("N5").Value = -5634.68

$ws.Range("H74").Value = 11500
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 45000
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = -47122

$ws.Range("H77").Value = 11500
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 135000
$ws.Range("M77").Value = -3696
$ws.Range("N77").Value = -145608

$ws.Range("H113").Value = 603
$ws.Range("I113").Value = 559.36365
$ws.Range("J113").Value = 731
$ws.Range("K113").Value = 1678.09095
$ws.Range("L113").Value = 2193
$ws.Range("M113").Value = 491.90905
$ws.Range("N113").Value = -6533

$ws.Range("H131").Value = 757.8
$ws.Range("J131").Value = 771.4583
$ws.Range("L131").Value = 2314.3749
$ws.Range("N131").Value = -12394.3749

$ws.Range("H135").Value = 1236.0513
$ws.Range("I135").Value = 222.64285
$ws.Range("J135").Value = 1803.56
$ws.Range("K135").Value = 2003.78565
$ws.Range("L135").Value = 16232.04
$ws.Range("M135").Value = 531.21435
$ws.Range("N135").Value = -21302.04

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 735.4138
$ws.Range("I97").Value = 673.88
$ws.Range("J97").Value = 1120
$ws.Range("K97").Value = 673.88
$ws.Range("L97").Value = 1120
$ws.Range("M97").Value = -177.88
$ws.Range("N97").Value = -2112

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1320
$ws.Range("I46").Value = 530.6
$ws.Range("J46").Value = 1491.6086
$ws.Range("K46").Value = 530.6
$ws.Range("L46").Value = 1491.6086
$ws.Range("M46").Value = -342.6
$ws.Range("N46").Value = -1867.6086

$ws.Range("H61").Value = 30304548
$ws.Range("I61").Value = 1871.4286
$ws.Range("J61").Value = 83334230
$ws.Range("K61").Value = 1871.4286
$ws.Range("L61").Value = 83334230
$ws.Range("M61").Value = -1669.4286
$ws.Range("N61").Value = -83334634

$ws.Range("H113").Value = 30304548
$ws.Range("I113").Value = 1871.4286
$ws.Range("J113").Value = 83334230
$ws.Range("K113").Value = 1871.4286
$ws.Range("L113").Value = 83334230
$ws.Range("M113").Value = 298.5714
$ws.Range("N113").Value = -83338570

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1887.3617
$ws.Range("I122").Value = 1571.7097
$ws.Range("J122").Value = 2498.9375
$ws.Range("K122").Value = 4715.1291
$ws.Range("L122").Value = 7496.8125
$ws.Range("M122").Value = -2265.1291
$ws.Range("N122").Value = -12396.8125

$ws.Range("H132").Value = 1474.7
$ws.Range("I132").Value = 1146.8
$ws.Range("J132").Value = 2786.3
$ws.Range("K132").Value = 3440.4
$ws.Range("L132").Value = 8358.900000000001
$ws.Range("M132").Value = -910.3999999999996
$ws.Range("N132").Value = -13418.9

$ws.Range("H133").Value = 46650.5
$ws.Range("J133").Value = 46650.5
$ws.Range("L133").Value = 46650.5
$ws.Range("N133").Value = -56770.5

$ws.Range("H135").Value = 82779.8
$ws.Range("J135").Value = 82779.8
$ws.Range("L135").Value = 82779.8
$ws.Range("N135").Value = -92919.8
